# Workbook: eex_ddh_JSON_lookup_complex.xlsx
# "added some mapping functions"
#
# The $-prefixed placeholder field names in column A (A2:A6) are renamed to
# plain mapping-function names (no leading "$"), the worksheet selection is
# moved off the old D2:D6 block, and column B is widened to fit its
# (longer) contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the leading "$" from the placeholder names in column A.
$ws.Range("A2").Value = "topic"
$ws.Range("A3").Value = "group"
$ws.Range("A4").Value = "release_date"
$ws.Range("A5").Value = "license_id"
$ws.Range("A6").Value = "format"

# Column B holds the longest strings (e.g. "field_external_metadata") -
# widen it to fit.
$ws.Columns("B").AutoFit() | Out-Null

# Move the selection.
$ws.Range("G8").Select() | Out-Null
